$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lithuania A Lyga")

# --- Re-order existing rows (tie-break shuffle from upstream re-sort) ---
# Row 89
$ws.Range("B89").Value = 6732827
$ws.Range("C89").Value = 'Lithuania A Lyga'
$ws.Range("D89").Value = 45220.375
$ws.Range("E89").Value = 'FK Dziugas Telsiai'
$ws.Range("F89").Value = 'FK Kauno Zalgiris'
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 2
$ws.Range("I89").Value = 'A'
$ws.Range("J89").Value = 6
$ws.Range("K89").Value = 3.9
$ws.Range("L89").Value = 1.444
$ws.Range("M89").Value = 4.75
$ws.Range("N89").Value = 3.6
$ws.Range("O89").Value = 1.65
$ws.Range("P89").Value = 0.75
$ws.Range("Q89").Value = 1.9
$ws.Range("R89").Value = 1.9
$ws.Range("S89").Value = 2.5
$ws.Range("T89").Value = 1.95
$ws.Range("U89").Value = 1.85
$ws.Range("V89").Value = -1
$ws.Range("W89").Value = -1
$ws.Range("X89").Value = 0.6499999999999999
$ws.Range("Y89").Value = -1
$ws.Range("Z89").Value = 0.8999999999999999
$ws.Range("AA89").Value = -1
$ws.Range("AB89").Value = 0.8500000000000001

# Row 90
$ws.Range("B90").Value = 7326568
$ws.Range("C90").Value = 'Lithuania A Lyga'
$ws.Range("D90").Value = 45220.375
$ws.Range("E90").Value = 'Hegelmann Litauen'
$ws.Range("F90").Value = 'Panevezys'
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 'D'
$ws.Range("J90").Value = 2.375
$ws.Range("K90").Value = 3.2
$ws.Range("L90").Value = 2.625
$ws.Range("M90").Value = 2.7
$ws.Range("N90").Value = 3.2
$ws.Range("O90").Value = 2.3
$ws.Range("P90").Value = 0
$ws.Range("Q90").Value = 2.05
$ws.Range("R90").Value = 1.75
$ws.Range("S90").Value = 2.25
$ws.Range("T90").Value = 1.875
$ws.Range("U90").Value = 1.925
$ws.Range("V90").Value = -1
$ws.Range("W90").Value = 2.2
$ws.Range("X90").Value = -1
$ws.Range("Y90").Value = 0
$ws.Range("Z90").Value = 0
$ws.Range("AA90").Value = -1
$ws.Range("AB90").Value = 0.925

# Row 100
$ws.Range("B100").Value = 6732727
$ws.Range("C100").Value = 'Lithuania A Lyga'
$ws.Range("D100").Value = 45242.41319444445
$ws.Range("E100").Value = 'FK Zalgiris Vilnius'
$ws.Range("F100").Value = 'FK Dainava Alytus'
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 'H'
$ws.Range("J100").Value = 1.285
$ws.Range("K100").Value = 5.5
$ws.Range("L100").Value = 6.5
$ws.Range("M100").Value = 1.3
$ws.Range("N100").Value = 5.5
$ws.Range("O100").Value = 6
$ws.Range("P100").Value = -1.5
$ws.Range("Q100").Value = 1.9
$ws.Range("R100").Value = 1.9
$ws.Range("S100").Value = 2.75
$ws.Range("T100").Value = 1.8
$ws.Range("U100").Value = 2
$ws.Range("V100").Value = 0.3
$ws.Range("W100").Value = -1
$ws.Range("X100").Value = -1
$ws.Range("Y100").Value = -1
$ws.Range("Z100").Value = 0.8999999999999999
$ws.Range("AA100").Value = -1
$ws.Range("AB100").Value = 1

# Row 101
$ws.Range("B101").Value = 6732836
$ws.Range("C101").Value = 'Lithuania A Lyga'
$ws.Range("D101").Value = 45242.41319444445
$ws.Range("E101").Value = 'FK Siauliai'
$ws.Range("F101").Value = 'Banga Gargzdai'
$ws.Range("G101").Value = 3
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 'H'
$ws.Range("J101").Value = 1.222
$ws.Range("K101").Value = 5.5
$ws.Range("L101").Value = 9
$ws.Range("M101").Value = 1.363
$ws.Range("N101").Value = 4.5
$ws.Range("O101").Value = 7
$ws.Range("P101").Value = -1.25
$ws.Range("Q101").Value = 1.9
$ws.Range("R101").Value = 1.9
$ws.Range("S101").Value = 2.5
$ws.Range("T101").Value = 1.975
$ws.Range("U101").Value = 1.825
$ws.Range("V101").Value = 0.363
$ws.Range("W101").Value = -1
$ws.Range("X101").Value = -1
$ws.Range("Y101").Value = 0.8999999999999999
$ws.Range("Z101").Value = -1
$ws.Range("AA101").Value = 0.9750000000000001
$ws.Range("AB101").Value = -1

# Row 102
$ws.Range("B102").Value = 6732837
$ws.Range("C102").Value = 'Lithuania A Lyga'
$ws.Range("D102").Value = 45242.41319444445
$ws.Range("E102").Value = 'Suduva Marijampole'
$ws.Range("F102").Value = 'FK Riteriai'
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 3
$ws.Range("I102").Value = 'A'
$ws.Range("J102").Value = 3.6
$ws.Range("K102").Value = 3.6
$ws.Range("L102").Value = 1.8
$ws.Range("M102").Value = 3
$ws.Range("N102").Value = 3.6
$ws.Range("O102").Value = 2
$ws.Range("P102").Value = 0.25
$ws.Range("Q102").Value = 2
$ws.Range("R102").Value = 1.8
$ws.Range("S102").Value = 2.5
$ws.Range("T102").Value = 1.975
$ws.Range("U102").Value = 1.825
$ws.Range("V102").Value = -1
$ws.Range("W102").Value = -1
$ws.Range("X102").Value = 1
$ws.Range("Y102").Value = -1
$ws.Range("Z102").Value = 0.8
$ws.Range("AA102").Value = 0.9750000000000001
$ws.Range("AB102").Value = -1

# Row 103
$ws.Range("B103").Value = 7465686
$ws.Range("C103").Value = 'Lithuania A Lyga'
$ws.Range("D103").Value = 45242.41319444445
$ws.Range("E103").Value = 'FK Kauno Zalgiris'
$ws.Range("F103").Value = 'Hegelmann Litauen'
$ws.Range("G103").Value = 4
$ws.Range("H103").Value = 2
$ws.Range("I103").Value = 'H'
$ws.Range("J103").Value = 2.3
$ws.Range("K103").Value = 4
$ws.Range("L103").Value = 2.3
$ws.Range("M103").Value = 2.55
$ws.Range("N103").Value = 4
$ws.Range("O103").Value = 2.2
$ws.Range("P103").Value = 0.25
$ws.Range("Q103").Value = 1.8
$ws.Range("R103").Value = 2
$ws.Range("S103").Value = 2.75
$ws.Range("T103").Value = 1.85
$ws.Range("U103").Value = 1.95
$ws.Range("V103").Value = 1.55
$ws.Range("W103").Value = -1
$ws.Range("X103").Value = -1
$ws.Range("Y103").Value = 0.8
$ws.Range("Z103").Value = -1
$ws.Range("AA103").Value = 0.8500000000000001
$ws.Range("AB103").Value = -1

# Row 104
$ws.Range("B104").Value = 6732834
$ws.Range("C104").Value = 'Lithuania A Lyga'
$ws.Range("D104").Value = 45242.41319444445
$ws.Range("E104").Value = 'Panevezys'
$ws.Range("F104").Value = 'FK Dziugas Telsiai'
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 'D'
$ws.Range("J104").Value = 1.25
$ws.Range("K104").Value = 5.5
$ws.Range("L104").Value = 7.5
$ws.Range("M104").Value = 1.45
$ws.Range("N104").Value = 4.5
$ws.Range("O104").Value = 5
$ws.Range("P104").Value = -1
$ws.Range("Q104").Value = 1.775
$ws.Range("R104").Value = 2.025
$ws.Range("S104").Value = 2.5
$ws.Range("T104").Value = 1.875
$ws.Range("U104").Value = 1.925
$ws.Range("V104").Value = -1
$ws.Range("W104").Value = 3.5
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 1.025
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = 0.925

# Row 136
$ws.Range("B136").Value = 7862044
$ws.Range("C136").Value = 'Lithuania A Lyga'
$ws.Range("D136").Value = 45392.5
$ws.Range("E136").Value = 'Banga Gargzdai'
$ws.Range("F136").Value = 'Suduva Marijampole'
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 'D'
$ws.Range("J136").Value = 2.1
$ws.Range("K136").Value = 3.05
$ws.Range("L136").Value = 3.2
$ws.Range("M136").Value = 3.4
$ws.Range("N136").Value = 3
$ws.Range("O136").Value = 2.1
$ws.Range("P136").Value = 0.25
$ws.Range("Q136").Value = 1.95
$ws.Range("R136").Value = 1.85
$ws.Range("S136").Value = 2
$ws.Range("T136").Value = 1.95
$ws.Range("U136").Value = 1.85
$ws.Range("V136").Value = -1
$ws.Range("W136").Value = 2
$ws.Range("X136").Value = -1
$ws.Range("Y136").Value = 0.475
$ws.Range("Z136").Value = -0.5
$ws.Range("AA136").Value = -1
$ws.Range("AB136").Value = 0.8500000000000001

# Row 137
$ws.Range("B137").Value = 7862922
$ws.Range("C137").Value = 'Lithuania A Lyga'
$ws.Range("D137").Value = 45392.5
$ws.Range("E137").Value = 'FK Siauliai'
$ws.Range("F137").Value = 'Panevezys'
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 1
$ws.Range("I137").Value = 'D'
$ws.Range("J137").Value = 2.7
$ws.Range("K137").Value = 3
$ws.Range("L137").Value = 2.5
$ws.Range("M137").Value = 2.9
$ws.Range("N137").Value = 2.9
$ws.Range("O137").Value = 2.375
$ws.Range("P137").Value = 0.25
$ws.Range("Q137").Value = 1.75
$ws.Range("R137").Value = 2.05
$ws.Range("S137").Value = 1.75
$ws.Range("T137").Value = 1.775
$ws.Range("U137").Value = 2.025
$ws.Range("V137").Value = -1
$ws.Range("W137").Value = 1.9
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = 0.375
$ws.Range("Z137").Value = -0.5
$ws.Range("AA137").Value = 0.3875
$ws.Range("AB137").Value = -0.5

# --- Append new match rows 163-167 ---
# Row 163
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A163").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D163").PasteSpecial(-4122) | Out-Null
$ws.Range("A163").Value = 161
$ws.Range("B163").Value = 7862940
$ws.Range("C163").Value = 'Lithuania A Lyga'
$ws.Range("D163").Value = 45423.375
$ws.Range("E163").Value = 'FK Dainava Alytus'
$ws.Range("F163").Value = 'FK Transinvest'
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 0
$ws.Range("I163").Value = 'H'
$ws.Range("J163").Value = 2.375
$ws.Range("K163").Value = 3
$ws.Range("L163").Value = 2.875
$ws.Range("M163").Value = 2.2
$ws.Range("N163").Value = 2.9
$ws.Range("O163").Value = 3.4
$ws.Range("P163").Value = -0.25
$ws.Range("Q163").Value = 1.925
$ws.Range("R163").Value = 1.875
$ws.Range("S163").Value = 2
$ws.Range("T163").Value = 1.9
$ws.Range("U163").Value = 1.9
$ws.Range("V163").Value = 1.2
$ws.Range("W163").Value = -1
$ws.Range("X163").Value = -1
$ws.Range("Y163").Value = 0.925
$ws.Range("Z163").Value = -1
$ws.Range("AA163").Value = -1
$ws.Range("AB163").Value = 0.8999999999999999

# Row 164
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A164").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D164").PasteSpecial(-4122) | Out-Null
$ws.Range("A164").Value = 162
$ws.Range("B164").Value = 7862054
$ws.Range("C164").Value = 'Lithuania A Lyga'
$ws.Range("D164").Value = 45423.45833333334
$ws.Range("E164").Value = 'Suduva Marijampole'
$ws.Range("F164").Value = 'FK Dziugas Telsiai'
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 1
$ws.Range("I164").Value = 'A'
$ws.Range("J164").Value = 2.15
$ws.Range("K164").Value = 3
$ws.Range("L164").Value = 3.25
$ws.Range("M164").Value = 2.05
$ws.Range("N164").Value = 2.875
$ws.Range("O164").Value = 3.75
$ws.Range("P164").Value = -0.25
$ws.Range("Q164").Value = 1.825
$ws.Range("R164").Value = 1.975
$ws.Range("S164").Value = 1.75
$ws.Range("T164").Value = 1.825
$ws.Range("U164").Value = 1.975
$ws.Range("V164").Value = -1
$ws.Range("W164").Value = -1
$ws.Range("X164").Value = 2.75
$ws.Range("Y164").Value = -1
$ws.Range("Z164").Value = 0.9750000000000001
$ws.Range("AA164").Value = -1
$ws.Range("AB164").Value = 0.9750000000000001

# Row 165
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A165").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D165").PasteSpecial(-4122) | Out-Null
$ws.Range("A165").Value = 163
$ws.Range("B165").Value = 7862941
$ws.Range("C165").Value = 'Lithuania A Lyga'
$ws.Range("D165").Value = 45423.54166666666
$ws.Range("E165").Value = 'Banga Gargzdai'
$ws.Range("F165").Value = 'Panevezys'
$ws.Range("G165").Value = 2
$ws.Range("H165").Value = 0
$ws.Range("I165").Value = 'H'
$ws.Range("J165").Value = 4.25
$ws.Range("K165").Value = 3.2
$ws.Range("L165").Value = 1.8
$ws.Range("M165").Value = 3.8
$ws.Range("N165").Value = 3.2
$ws.Range("O165").Value = 1.85
$ws.Range("P165").Value = 0.5
$ws.Range("Q165").Value = 1.875
$ws.Range("R165").Value = 1.925
$ws.Range("S165").Value = 2
$ws.Range("T165").Value = 1.775
$ws.Range("U165").Value = 2.025
$ws.Range("V165").Value = 2.8
$ws.Range("W165").Value = -1
$ws.Range("X165").Value = -1
$ws.Range("Y165").Value = 0.875
$ws.Range("Z165").Value = -1
$ws.Range("AA165").Value = 0
$ws.Range("AB165").Value = 0

# Row 166
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A166").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D166").PasteSpecial(-4122) | Out-Null
$ws.Range("A166").Value = 164
$ws.Range("B166").Value = 7862055
$ws.Range("C166").Value = 'Lithuania A Lyga'
$ws.Range("D166").Value = 45424.41666666666
$ws.Range("E166").Value = 'Hegelmann Litauen'
$ws.Range("F166").Value = 'FK Zalgiris Vilnius'
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0
$ws.Range("I166").Value = 'D'
$ws.Range("J166").Value = 3.75
$ws.Range("K166").Value = 3.4
$ws.Range("L166").Value = 1.85
$ws.Range("M166").Value = 3.6
$ws.Range("N166").Value = 3.8
$ws.Range("O166").Value = 1.8
$ws.Range("P166").Value = 0.5
$ws.Range("Q166").Value = 1.975
$ws.Range("R166").Value = 1.825
$ws.Range("S166").Value = 3
$ws.Range("T166").Value = 1.95
$ws.Range("U166").Value = 1.85
$ws.Range("V166").Value = -1
$ws.Range("W166").Value = 2.8
$ws.Range("X166").Value = -1
$ws.Range("Y166").Value = 0.9750000000000001
$ws.Range("Z166").Value = -1
$ws.Range("AA166").Value = -1
$ws.Range("AB166").Value = 0.8500000000000001

# Row 167
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A167").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D167").PasteSpecial(-4122) | Out-Null
$ws.Range("A167").Value = 165
$ws.Range("B167").Value = 7862942
$ws.Range("C167").Value = 'Lithuania A Lyga'
$ws.Range("D167").Value = 45424.51736111111
$ws.Range("E167").Value = 'FK Kauno Zalgiris'
$ws.Range("F167").Value = 'FK Siauliai'
$ws.Range("G167").Value = 2
$ws.Range("H167").Value = 1
$ws.Range("I167").Value = 'H'
$ws.Range("J167").Value = 2.05
$ws.Range("K167").Value = 3.1
$ws.Range("L167").Value = 3.4
$ws.Range("M167").Value = 2.05
$ws.Range("N167").Value = 3.2
$ws.Range("O167").Value = 3.5
$ws.Range("P167").Value = -0.25
$ws.Range("Q167").Value = 1.775
$ws.Range("R167").Value = 2.025
$ws.Range("S167").Value = 2.5
$ws.Range("T167").Value = 2
$ws.Range("U167").Value = 1.8
$ws.Range("V167").Value = 1.05
$ws.Range("W167").Value = -1
$ws.Range("X167").Value = -1
$ws.Range("Y167").Value = 0.7749999999999999
$ws.Range("Z167").Value = -1
$ws.Range("AA167").Value = 1
$ws.Range("AB167").Value = -1
